$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "title"
$ws.Range("B3").Value = "Food Prep"

$ws.Range("A4").Value = "options"
$ws.Range("B4").Value = "OPTIONS"

$ws.Range("A5").Value = "music"
$ws.Range("B5").Value = "MUSIC"

$ws.Range("A6").Value = "sound"
$ws.Range("B6").Value = "SOUND"

$ws.Range("A7").Value = "close"
$ws.Range("B7").Value = "CLOSE"

$ws.Range("A8").Select()
